$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.304.82"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "'3.106.57"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'525.51"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").Value = "'138.01"
$ws.Range("E6").Value = "  -2.83%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'3.107.63"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").Value = "'0.452"
$ws.Range("E9").Value = "  +2.35%  "
$ws.Range("D10").Value = "'7.35"
$ws.Range("E10").Value = "  +1.87%  "
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("D12").Value = "'0.402"
$ws.Range("E12").Value = "  +2.50%  "
$ws.Range("D13").Value = "'3.637.63"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("E14").Value = "  +1.74%  "
$ws.Range("D15").Value = "'25.58"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("D17").Value = "'57.445.09"
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").Value = "'3.101.18"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("E19").Value = "  -2.66%  "
$ws.Range("D20").Value = "'12.55"
$ws.Range("E20").Value = "  -1.65%  "
$ws.Range("E21").Value = "  -0.75%  "
$ws.Range("D22").Value = "'350.16"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("D24").Value = "'68.40"
$ws.Range("E24").Value = "  +1.73%  "
$ws.Range("D25").Value = "'0.503"
$ws.Range("E25").Value = "  -1.84%  "
$ws.Range("D26").Value = "'0.168"
$ws.Range("E26").Value = "  -0.92%  "
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").Value = "'0.0₃0891"
$ws.Range("E28").Value = "  -2.75%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "'7.37"
$ws.Range("E30").Value = "  +1.49%  "
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("D32").Value = "'5.99"
$ws.Range("E32").Value = "  -7.25%  "
$ws.Range("D33").Value = "'20.88"
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("E34").Value = "  +8.10%  "
$ws.Range("E35").Value = "  -3.74%  "
$ws.Range("D36").Value = "'159.49"
$ws.Range("E36").Value = "  +0.77%  "
$ws.Range("E37").Value = "  -1.58%  "
$ws.Range("D38").Value = "'26.43"
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("E40").Value = "  -1.37%  "
$ws.Range("E42").Value = "  +1.64%  "
$ws.Range("D43").Value = "'0.697"
$ws.Range("E43").Value = "  +2.01%  "
$ws.Range("D44").Value = "'2.408.54"
$ws.Range("E44").Value = "  +5.34%  "
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").Value = "'3.146.02"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("D49").Value = "'0.973"
$ws.Range("E49").Value = "  -2.40%  "
$ws.Range("D50").Value = "'5.99"
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'19.79"
$ws.Range("E51").Value = "  -4.09%  "
